$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows with new/changed values ---
# Row 114
$ws.Cells.Item(114, 2).Value = 7559468
$ws.Cells.Item(114, 6).Value = "Liverpool Montevideo"
$ws.Cells.Item(114, 7).Value = "CA River Plate"
$ws.Cells.Item(114, 8).Value = 2
$ws.Cells.Item(114, 9).Value = 1
$ws.Cells.Item(114, 10).Value = "H"
$ws.Cells.Item(114, 11).Value = 1.7
$ws.Cells.Item(114, 12).Value = 3
$ws.Cells.Item(114, 13).Value = 5.75
$ws.Cells.Item(114, 14).Value = 1.833
$ws.Cells.Item(114, 16).Value = 4.5
$ws.Cells.Item(114, 17).Value = -0.5
$ws.Cells.Item(114, 18).Value = 1.925
$ws.Cells.Item(114, 19).Value = 1.925
$ws.Cells.Item(114, 20).Value = 2.25
$ws.Cells.Item(114, 21).Value = 2.025
$ws.Cells.Item(114, 22).Value = 1.825
$ws.Cells.Item(114, 23).Value = 0.833
$ws.Cells.Item(114, 24).Value = -1
$ws.Cells.Item(114, 26).Value = 0.925
$ws.Cells.Item(114, 27).Value = -1
$ws.Cells.Item(114, 28).Value = 1.025
$ws.Cells.Item(114, 29).Value = -1

# Row 115
$ws.Cells.Item(115, 2).Value = 7559469
$ws.Cells.Item(115, 6).Value = "Montevideo Wanderers"
$ws.Cells.Item(115, 7).Value = "Penarol"
$ws.Cells.Item(115, 8).Value = 0
$ws.Cells.Item(115, 9).Value = 0
$ws.Cells.Item(115, 10).Value = "D"
$ws.Cells.Item(115, 11).Value = 4.75
$ws.Cells.Item(115, 12).Value = 3.4
$ws.Cells.Item(115, 13).Value = 1.7
$ws.Cells.Item(115, 14).Value = 2.7
$ws.Cells.Item(115, 16).Value = 2.45
$ws.Cells.Item(115, 17).Value = 0
$ws.Cells.Item(115, 18).Value = 2.05
$ws.Cells.Item(115, 19).Value = 1.8
$ws.Cells.Item(115, 20).Value = 2.5
$ws.Cells.Item(115, 21).Value = 1.975
$ws.Cells.Item(115, 22).Value = 1.875
$ws.Cells.Item(115, 23).Value = -1
$ws.Cells.Item(115, 24).Value = 2.2
$ws.Cells.Item(115, 26).Value = 0
$ws.Cells.Item(115, 27).Value = -0
$ws.Cells.Item(115, 28).Value = -1
$ws.Cells.Item(115, 29).Value = 0.875

# Row 117
$ws.Cells.Item(117, 2).Value = 7013885
$ws.Cells.Item(117, 6).Value = "La Luz"
$ws.Cells.Item(117, 7).Value = "Atletico Fenix Montevideo"
$ws.Cells.Item(117, 8).Value = 0
$ws.Cells.Item(117, 9).Value = 2
$ws.Cells.Item(117, 10).Value = "A"
$ws.Cells.Item(117, 11).Value = 3
$ws.Cells.Item(117, 12).Value = 3
$ws.Cells.Item(117, 13).Value = 2.4
$ws.Cells.Item(117, 14).Value = 2.9
$ws.Cells.Item(117, 15).Value = 2.75
$ws.Cells.Item(117, 16).Value = 2.6
$ws.Cells.Item(117, 17).Value = 0
$ws.Cells.Item(117, 18).Value = 2.025
$ws.Cells.Item(117, 19).Value = 1.825
$ws.Cells.Item(117, 20).Value = 2
$ws.Cells.Item(117, 21).Value = 2.025
$ws.Cells.Item(117, 22).Value = 1.825
$ws.Cells.Item(117, 24).Value = -1
$ws.Cells.Item(117, 25).Value = 1.6
$ws.Cells.Item(117, 27).Value = 0.825
$ws.Cells.Item(117, 28).Value = 0
$ws.Cells.Item(117, 29).Value = -0

# Row 119
$ws.Cells.Item(119, 2).Value = 7013702
$ws.Cells.Item(119, 6).Value = "Defensor Sporting"
$ws.Cells.Item(119, 7).Value = "Danubio"
$ws.Cells.Item(119, 11).Value = 1.8
$ws.Cells.Item(119, 12).Value = 3.6
$ws.Cells.Item(119, 13).Value = 4.2
$ws.Cells.Item(119, 14).Value = 1.8
$ws.Cells.Item(119, 15).Value = 3.6
$ws.Cells.Item(119, 16).Value = 4.2
$ws.Cells.Item(119, 17).Value = -0.75
$ws.Cells.Item(119, 18).Value = 2.05
$ws.Cells.Item(119, 19).Value = 1.8
$ws.Cells.Item(119, 20).Value = 2.25
$ws.Cells.Item(119, 21).Value = 1.85
$ws.Cells.Item(119, 22).Value = 2
$ws.Cells.Item(119, 25).Value = 3.2
$ws.Cells.Item(119, 27).Value = 0.8
$ws.Cells.Item(119, 28).Value = -0.5
$ws.Cells.Item(119, 29).Value = 0.5

# Row 120
$ws.Cells.Item(120, 2).Value = 7013409
$ws.Cells.Item(120, 6).Value = "Nacional De Football"
$ws.Cells.Item(120, 7).Value = "Torque"
$ws.Cells.Item(120, 8).Value = 1
$ws.Cells.Item(120, 9).Value = 1
$ws.Cells.Item(120, 10).Value = "D"
$ws.Cells.Item(120, 11).Value = 1.666
$ws.Cells.Item(120, 12).Value = 3.9
$ws.Cells.Item(120, 13).Value = 4.5
$ws.Cells.Item(120, 14).Value = 1.615
$ws.Cells.Item(120, 15).Value = 4
$ws.Cells.Item(120, 16).Value = 4.75
$ws.Cells.Item(120, 18).Value = 1.8
$ws.Cells.Item(120, 19).Value = 2.05
$ws.Cells.Item(120, 20).Value = 2.75
$ws.Cells.Item(120, 21).Value = 1.95
$ws.Cells.Item(120, 22).Value = 1.9
$ws.Cells.Item(120, 24).Value = 3
$ws.Cells.Item(120, 25).Value = -1
$ws.Cells.Item(120, 27).Value = 1.05
$ws.Cells.Item(120, 28).Value = -1
$ws.Cells.Item(120, 29).Value = 0.8999999999999999

# Row 169
$ws.Cells.Item(169, 2).Value = 8014131
$ws.Cells.Item(169, 5).Value = 45387.79166666666
$ws.Cells.Item(169, 6).Value = "Boston River"
$ws.Cells.Item(169, 7).Value = "Defensor Sporting"
$ws.Cells.Item(169, 8).Value = 2
$ws.Cells.Item(169, 9).Value = 2
$ws.Cells.Item(169, 10).Value = "D"
$ws.Cells.Item(169, 11).Value = 2.875
$ws.Cells.Item(169, 12).Value = 3.3
$ws.Cells.Item(169, 13).Value = 2.375
$ws.Cells.Item(169, 14).Value = 3.1
$ws.Cells.Item(169, 15).Value = 3.3
$ws.Cells.Item(169, 16).Value = 2.2
$ws.Cells.Item(169, 17).Value = 0.25
$ws.Cells.Item(169, 18).Value = 1.9
$ws.Cells.Item(169, 19).Value = 1.95
$ws.Cells.Item(169, 21).Value = 1.825
$ws.Cells.Item(169, 22).Value = 2.025
$ws.Cells.Item(169, 23).Value = -1
$ws.Cells.Item(169, 24).Value = 2.3
$ws.Cells.Item(169, 25).Value = -1
$ws.Cells.Item(169, 26).Value = 0.45
$ws.Cells.Item(169, 27).Value = -0.5
$ws.Cells.Item(169, 28).Value = 0.825
$ws.Cells.Item(169, 29).Value = -1

# Row 170
$ws.Cells.Item(170, 2).Value = 8014133
$ws.Cells.Item(170, 5).Value = 45388.41666666666
$ws.Cells.Item(170, 6).Value = "CA River Plate"
$ws.Cells.Item(170, 7).Value = "Montevideo Wanderers"
$ws.Cells.Item(170, 8).Value = 3
$ws.Cells.Item(170, 9).Value = 1
$ws.Cells.Item(170, 10).Value = "H"
$ws.Cells.Item(170, 11).Value = 2.5
$ws.Cells.Item(170, 12).Value = 2.9
$ws.Cells.Item(170, 13).Value = 2.9
$ws.Cells.Item(170, 14).Value = 2.7
$ws.Cells.Item(170, 15).Value = 2.9
$ws.Cells.Item(170, 16).Value = 2.7
$ws.Cells.Item(170, 17).Value = 0
$ws.Cells.Item(170, 18).Value = 1.95
$ws.Cells.Item(170, 19).Value = 1.9
$ws.Cells.Item(170, 20).Value = 2
$ws.Cells.Item(170, 21).Value = 1.875
$ws.Cells.Item(170, 22).Value = 1.975
$ws.Cells.Item(170, 23).Value = 1.7
$ws.Cells.Item(170, 24).Value = -1
$ws.Cells.Item(170, 25).Value = -1
$ws.Cells.Item(170, 26).Value = 0.95
$ws.Cells.Item(170, 27).Value = -1
$ws.Cells.Item(170, 28).Value = 0.875
$ws.Cells.Item(170, 29).Value = -1

# Row 171
$ws.Cells.Item(171, 2).Value = 8014044
$ws.Cells.Item(171, 5).Value = 45388.52083333334
$ws.Cells.Item(171, 6).Value = "Racing Club de Montevideo"
$ws.Cells.Item(171, 7).Value = "Miramar Misiones"
$ws.Cells.Item(171, 8).Value = 1
$ws.Cells.Item(171, 9).Value = 1
$ws.Cells.Item(171, 10).Value = "D"
$ws.Cells.Item(171, 11).Value = 1.909
$ws.Cells.Item(171, 12).Value = 3.3
$ws.Cells.Item(171, 13).Value = 4
$ws.Cells.Item(171, 14).Value = 1.909
$ws.Cells.Item(171, 15).Value = 3.4
$ws.Cells.Item(171, 16).Value = 4
$ws.Cells.Item(171, 17).Value = -0.5
$ws.Cells.Item(171, 18).Value = 1.875
$ws.Cells.Item(171, 19).Value = 1.975
$ws.Cells.Item(171, 20).Value = 2.25
$ws.Cells.Item(171, 21).Value = 1.8
$ws.Cells.Item(171, 22).Value = 2.05
$ws.Cells.Item(171, 23).Value = -1
$ws.Cells.Item(171, 24).Value = 2.4
$ws.Cells.Item(171, 25).Value = -1
$ws.Cells.Item(171, 26).Value = -1
$ws.Cells.Item(171, 27).Value = 0.9750000000000001
$ws.Cells.Item(171, 28).Value = -0.5
$ws.Cells.Item(171, 29).Value = 0.5249999999999999

# Row 172
$ws.Cells.Item(172, 2).Value = 8014043
$ws.Cells.Item(172, 5).Value = 45388.625
$ws.Cells.Item(172, 6).Value = "Danubio"
$ws.Cells.Item(172, 7).Value = "Club Atletico Progreso"
$ws.Cells.Item(172, 8).Value = 1
$ws.Cells.Item(172, 9).Value = 1
$ws.Cells.Item(172, 10).Value = "D"
$ws.Cells.Item(172, 11).Value = 2.375
$ws.Cells.Item(172, 12).Value = 3.1
$ws.Cells.Item(172, 13).Value = 3
$ws.Cells.Item(172, 15).Value = 3.25
$ws.Cells.Item(172, 20).Value = 2.5
$ws.Cells.Item(172, 21).Value = 2.025
$ws.Cells.Item(172, 22).Value = 1.825
$ws.Cells.Item(172, 23).Value = -1
$ws.Cells.Item(172, 24).Value = 2.25
$ws.Cells.Item(172, 25).Value = -1
$ws.Cells.Item(172, 26).Value = -0.5
$ws.Cells.Item(172, 27).Value = 0.45
$ws.Cells.Item(172, 28).Value = -1
$ws.Cells.Item(172, 29).Value = 0.825

# Row 173
$ws.Cells.Item(173, 2).Value = 8014089
$ws.Cells.Item(173, 5).Value = 45388.75
$ws.Cells.Item(173, 6).Value = "Nacional De Football"
$ws.Cells.Item(173, 7).Value = "Cerro Largo"
$ws.Cells.Item(173, 8).Value = 1
$ws.Cells.Item(173, 9).Value = 0
$ws.Cells.Item(173, 10).Value = "H"
$ws.Cells.Item(173, 11).Value = 1.5
$ws.Cells.Item(173, 12).Value = 4
$ws.Cells.Item(173, 13).Value = 6
$ws.Cells.Item(173, 14).Value = 1.6
$ws.Cells.Item(173, 15).Value = 3.8
$ws.Cells.Item(173, 16).Value = 5
$ws.Cells.Item(173, 17).Value = -0.75
$ws.Cells.Item(173, 18).Value = 1.8
$ws.Cells.Item(173, 19).Value = 2.05
$ws.Cells.Item(173, 20).Value = 2.25
$ws.Cells.Item(173, 21).Value = 1.85
$ws.Cells.Item(173, 22).Value = 2
$ws.Cells.Item(173, 23).Value = 0.6000000000000001
$ws.Cells.Item(173, 24).Value = -1
$ws.Cells.Item(173, 25).Value = -1
$ws.Cells.Item(173, 26).Value = 0.4
$ws.Cells.Item(173, 27).Value = -0.5
$ws.Cells.Item(173, 28).Value = -1
$ws.Cells.Item(173, 29).Value = 1

# --- Add new rows, copying formatting from the last existing row (173) ---
$ws.Range("A173:AC173").Copy()
$ws.Range("A174:AC174").PasteSpecial(-4122)
$ws.Range("A173:AC173").Copy()
$ws.Range("A175:AC175").PasteSpecial(-4122)
$ws.Range("A173:AC173").Copy()
$ws.Range("A176:AC176").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 174 (new)
$ws.Cells.Item(174, 1).Value = 172
$ws.Cells.Item(174, 2).Value = 8014132
$ws.Cells.Item(174, 3).Value = "Uruguay Primera División"
$ws.Cells.Item(174, 4).Value = "Uruguay Apertura"
$ws.Cells.Item(174, 5).Value = 45389.41666666666
$ws.Cells.Item(174, 6).Value = "Cerro"
$ws.Cells.Item(174, 7).Value = "Rampla Juniors"
$ws.Cells.Item(174, 8).Value = 3
$ws.Cells.Item(174, 9).Value = 0
$ws.Cells.Item(174, 10).Value = "H"
$ws.Cells.Item(174, 11).Value = 2.2
$ws.Cells.Item(174, 12).Value = 3.2
$ws.Cells.Item(174, 13).Value = 3.2
$ws.Cells.Item(174, 14).Value = 2.2
$ws.Cells.Item(174, 15).Value = 3.2
$ws.Cells.Item(174, 16).Value = 3.25
$ws.Cells.Item(174, 17).Value = -0.25
$ws.Cells.Item(174, 18).Value = 1.875
$ws.Cells.Item(174, 19).Value = 1.975
$ws.Cells.Item(174, 20).Value = 2.25
$ws.Cells.Item(174, 21).Value = 1.95
$ws.Cells.Item(174, 22).Value = 1.9
$ws.Cells.Item(174, 23).Value = 1.2
$ws.Cells.Item(174, 24).Value = -1
$ws.Cells.Item(174, 25).Value = -1
$ws.Cells.Item(174, 26).Value = 0.875
$ws.Cells.Item(174, 27).Value = -1
$ws.Cells.Item(174, 28).Value = 0.95
$ws.Cells.Item(174, 29).Value = -1

# Row 175 (new)
$ws.Cells.Item(175, 1).Value = 173
$ws.Cells.Item(175, 2).Value = 8014090
$ws.Cells.Item(175, 3).Value = "Uruguay Primera División"
$ws.Cells.Item(175, 4).Value = "Uruguay Apertura"
$ws.Cells.Item(175, 5).Value = 45389.625
$ws.Cells.Item(175, 6).Value = "Atletico Fenix Montevideo"
$ws.Cells.Item(175, 7).Value = "Liverpool Montevideo"
$ws.Cells.Item(175, 8).Value = 1
$ws.Cells.Item(175, 9).Value = 3
$ws.Cells.Item(175, 10).Value = "A"
$ws.Cells.Item(175, 11).Value = 2.625
$ws.Cells.Item(175, 12).Value = 2.9
$ws.Cells.Item(175, 13).Value = 2.9
$ws.Cells.Item(175, 14).Value = 2.55
$ws.Cells.Item(175, 15).Value = 2.9
$ws.Cells.Item(175, 16).Value = 2.9
$ws.Cells.Item(175, 17).Value = 0
$ws.Cells.Item(175, 18).Value = 1.825
$ws.Cells.Item(175, 19).Value = 2.025
$ws.Cells.Item(175, 20).Value = 2
$ws.Cells.Item(175, 21).Value = 1.925
$ws.Cells.Item(175, 22).Value = 1.925
$ws.Cells.Item(175, 23).Value = -1
$ws.Cells.Item(175, 24).Value = -1
$ws.Cells.Item(175, 25).Value = 1.9
$ws.Cells.Item(175, 26).Value = -1
$ws.Cells.Item(175, 27).Value = 1.025
$ws.Cells.Item(175, 28).Value = 0.925
$ws.Cells.Item(175, 29).Value = -1

# Row 176 (new)
$ws.Cells.Item(176, 1).Value = 174
$ws.Cells.Item(176, 2).Value = 8014091
$ws.Cells.Item(176, 3).Value = "Uruguay Primera División"
$ws.Cells.Item(176, 4).Value = "Uruguay Apertura"
$ws.Cells.Item(176, 5).Value = 45389.75
$ws.Cells.Item(176, 6).Value = "Deportivo Maldonado"
$ws.Cells.Item(176, 7).Value = "Penarol"
$ws.Cells.Item(176, 8).Value = 1
$ws.Cells.Item(176, 9).Value = 2
$ws.Cells.Item(176, 10).Value = "A"
$ws.Cells.Item(176, 11).Value = 5
$ws.Cells.Item(176, 12).Value = 3.75
$ws.Cells.Item(176, 13).Value = 1.615
$ws.Cells.Item(176, 14).Value = 3.1
$ws.Cells.Item(176, 15).Value = 3.3
$ws.Cells.Item(176, 16).Value = 2.1
$ws.Cells.Item(176, 17).Value = 0.25
$ws.Cells.Item(176, 18).Value = 1.95
$ws.Cells.Item(176, 19).Value = 1.9
$ws.Cells.Item(176, 20).Value = 2.25
$ws.Cells.Item(176, 21).Value = 1.85
$ws.Cells.Item(176, 22).Value = 2
$ws.Cells.Item(176, 23).Value = -1
$ws.Cells.Item(176, 24).Value = -1
$ws.Cells.Item(176, 25).Value = 1.1
$ws.Cells.Item(176, 26).Value = -1
$ws.Cells.Item(176, 27).Value = 0.8999999999999999
$ws.Cells.Item(176, 28).Value = 0.8500000000000001
$ws.Cells.Item(176, 29).Value = -1

